# Auto-generated Excel COM-interop script
# Applies updated market-board derived values (currentAveragePrice* / LevePrice* / LeveProfit*)
# to the per-class Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 92.333336
$ws.Range("J4").Value = 89.5
$ws.Range("L4").Value = 89.5
$ws.Range("N4").Value = -317.5
$ws.Range("H18").Value = 248.33333
$ws.Range("I18").Value = 248.33333
$ws.Range("K18").Value = 248.33333
$ws.Range("M18").Value = 35.66667000000001
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H38").Value = 1042
$ws.Range("I38").Value = 63
$ws.Range("K38").Value = 189
$ws.Range("M38").Value = 183
$ws.Range("H41").Value = 2579.875
$ws.Range("I41").Value = 3207.7273
$ws.Range("K41").Value = 3207.7273
$ws.Range("M41").Value = -2767.7273
$ws.Range("H43").Value = 1214.4286
$ws.Range("J43").Value = 1231.6666
$ws.Range("L43").Value = 1231.6666
$ws.Range("N43").Value = -1369.6666
$ws.Range("H58").Value = 181.09091
$ws.Range("I58").Value = 181.09091
$ws.Range("K58").Value = 543.27273
$ws.Range("M58").Value = -393.27273
$ws.Range("H62").Value = 276998.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 276998.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 99040.45
$ws.Range("I74").Value = 116721.664
$ws.Range("K74").Value = 116721.664
$ws.Range("M74").Value = -115785.664
$ws.Range("H77").Value = 99040.45
$ws.Range("I77").Value = 116721.664
$ws.Range("K77").Value = 583608.3200000001
$ws.Range("M77").Value = -578928.3200000001
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 1305
$ws.Range("I92").Value = 457.5
$ws.Range("K92").Value = 457.5
$ws.Range("M92").Value = 790.5
$ws.Range("H98").Value = 7750
$ws.Range("I98").Value = 7500
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 7500
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = -6002
$ws.Range("N98").Value = -10996
$ws.Range("H100").Value = 2233
$ws.Range("I100").Value = 2233
$ws.Range("K100").Value = 2233
$ws.Range("M100").Value = -1692
$ws.Range("H106").Value = 3726.5715
$ws.Range("I106").Value = 3709.6667
$ws.Range("K106").Value = 3709.6667
$ws.Range("M106").Value = -3078.6667
$ws.Range("H113").Value = 78996.81
$ws.Range("I113").Value = 96162.414
$ws.Range("K113").Value = 96162.414
$ws.Range("M113").Value = -92908.414
$ws.Range("H122").Value = 7750
$ws.Range("I122").Value = 7500
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 22500
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -20050
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 5409634
$ws.Range("I132").Value = 6430447
$ws.Range("J132").Value = 101405.9
$ws.Range("K132").Value = 19291341
$ws.Range("L132").Value = 304217.7
$ws.Range("M132").Value = -19288811
$ws.Range("N132").Value = -309277.7
$ws.Range("H141").Value = 2831.2424
$ws.Range("I141").Value = 2908.6072
$ws.Range("J141").Value = 2398
$ws.Range("K141").Value = 8725.821599999999
$ws.Range("L141").Value = 7194
$ws.Range("M141").Value = -3545.821599999999
$ws.Range("N141").Value = -17554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 19979.4
$ws.Range("I55").Value = 12449.5
$ws.Range("K55").Value = 12449.5
$ws.Range("M55").Value = -12134.5
$ws.Range("H76").Value = 48999
$ws.Range("J76").Value = 48999
$ws.Range("L76").Value = 48999
$ws.Range("N76").Value = -49675
$ws.Range("H79").Value = 48999
$ws.Range("J79").Value = 48999
$ws.Range("L79").Value = 48999
$ws.Range("N79").Value = -51339
$ws.Range("H122").Value = 3225
$ws.Range("I122").Value = 2982.4736
$ws.Range("K122").Value = 8947.4208
$ws.Range("M122").Value = -6497.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 109990
$ws.Range("J126").Value = 109990
$ws.Range("L126").Value = 109990
$ws.Range("N126").Value = -119870
$ws.Range("H134").Value = 1728.1666
$ws.Range("I134").Value = 1210.3549
$ws.Range("J134").Value = 4938.6
$ws.Range("K134").Value = 3631.0647
$ws.Range("L134").Value = 14815.8
$ws.Range("M134").Value = -1096.0647
$ws.Range("N134").Value = -19885.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.07143
$ws.Range("I7").Value = 75.75
$ws.Range("J7").Value = 340.16666
$ws.Range("K7").Value = 75.75
$ws.Range("L7").Value = 340.16666
$ws.Range("M7").Value = 37.25
$ws.Range("N7").Value = -566.16666
$ws.Range("H16").Value = 1987.7
$ws.Range("I16").Value = 1541.8889
$ws.Range("K16").Value = 1541.8889
$ws.Range("M16").Value = -1254.8889
$ws.Range("H22").Value = 10002
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 20865896
$ws.Range("I31").Value = 2815019.2
$ws.Range("J31").Value = 111120280
$ws.Range("K31").Value = 2815019.2
$ws.Range("L31").Value = 111120280
$ws.Range("M31").Value = -2814724.2
$ws.Range("N31").Value = -111120870
$ws.Range("H34").Value = 20865896
$ws.Range("I34").Value = 2815019.2
$ws.Range("J34").Value = 111120280
$ws.Range("K34").Value = 2815019.2
$ws.Range("L34").Value = 111120280
$ws.Range("M34").Value = -2814817.2
$ws.Range("N34").Value = -111120684
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H50").Value = 56500
$ws.Range("J50").Value = 56500
$ws.Range("L50").Value = 56500
$ws.Range("N50").Value = -57750
$ws.Range("H59").Value = 26000
$ws.Range("J59").Value = 26000
$ws.Range("L59").Value = 26000
$ws.Range("N59").Value = -28290
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H99").Value = 3615
$ws.Range("I99").Value = 3593.1
$ws.Range("J99").Value = 3724.5
$ws.Range("K99").Value = 3593.1
$ws.Range("L99").Value = 3724.5
$ws.Range("M99").Value = -2095.1
$ws.Range("N99").Value = -6720.5
$ws.Range("H113").Value = 1987.7
$ws.Range("I113").Value = 1541.8889
$ws.Range("K113").Value = 1541.8889
$ws.Range("M113").Value = 628.1111000000001
$ws.Range("H118").Value = 203999.5
$ws.Range("J118").Value = 203999.5
$ws.Range("L118").Value = 203999.5
$ws.Range("N118").Value = -207313.5
$ws.Range("H126").Value = 3615
$ws.Range("I126").Value = 3593.1
$ws.Range("J126").Value = 3724.5
$ws.Range("K126").Value = 10779.3
$ws.Range("L126").Value = 11173.5
$ws.Range("M126").Value = -8309.299999999999
$ws.Range("N126").Value = -16113.5
$ws.Range("H140").Value = 75439.164
$ws.Range("J140").Value = 75439.164
$ws.Range("L140").Value = 75439.164
$ws.Range("N140").Value = -85799.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 311.34616
$ws.Range("I6").Value = 311.34616
$ws.Range("K6").Value = 934.0384799999999
$ws.Range("M6").Value = -821.0384799999999
$ws.Range("H10").Value = 333.08334
$ws.Range("I10").Value = 110.77778
$ws.Range("K10").Value = 332.33334
$ws.Range("M10").Value = -193.33334
$ws.Range("H11").Value = 142975.14
$ws.Range("J11").Value = 333382.34
$ws.Range("L11").Value = 1000147.02
$ws.Range("N11").Value = -1000427.02
$ws.Range("H13").Value = 19.666666
$ws.Range("I13").Value = 19.666666
$ws.Range("K13").Value = 58.999998
$ws.Range("M13").Value = 109.000002
$ws.Range("H16").Value = 66.333336
$ws.Range("I16").Value = 66.333336
$ws.Range("K16").Value = 199.000008
$ws.Range("M16").Value = -26.00000800000001
$ws.Range("H17").Value = 415.42856
$ws.Range("J17").Value = 20
$ws.Range("L17").Value = 60
$ws.Range("N17").Value = -398
$ws.Range("H25").Value = 99.75
$ws.Range("I25").Value = 99.75
$ws.Range("K25").Value = 299.25
$ws.Range("M25").Value = -130.25
$ws.Range("H30").Value = 99.75
$ws.Range("I30").Value = 99.75
$ws.Range("K30").Value = 299.25
$ws.Range("M30").Value = -197.25
$ws.Range("H34").Value = 3863695
$ws.Range("I34").Value = 842980
$ws.Range("J34").Value = 11113411
$ws.Range("K34").Value = 2528940
$ws.Range("L34").Value = 33340233
$ws.Range("M34").Value = -2528856
$ws.Range("N34").Value = -33340401
$ws.Range("H39").Value = 2335.875
$ws.Range("J39").Value = 2557.4
$ws.Range("L39").Value = 7672.200000000001
$ws.Range("N39").Value = -8260.200000000001
$ws.Range("H55").Value = 2901.875
$ws.Range("I55").Value = 226.66667
$ws.Range("J55").Value = 3519.2307
$ws.Range("K55").Value = 680.00001
$ws.Range("L55").Value = 10557.6921
$ws.Range("M55").Value = -503.00001
$ws.Range("N55").Value = -10911.6921
$ws.Range("H68").Value = 3572813.2
$ws.Range("J68").Value = 5001540.5
$ws.Range("L68").Value = 15004621.5
$ws.Range("N68").Value = -15006243.5
$ws.Range("H71").Value = 3572813.2
$ws.Range("J71").Value = 5001540.5
$ws.Range("L71").Value = 45013864.5
$ws.Range("N71").Value = -45021976.5
$ws.Range("H131").Value = 18042.766
$ws.Range("I131").Value = 112121.555
$ws.Range("J131").Value = 2648.0544
$ws.Range("K131").Value = 336364.665
$ws.Range("L131").Value = 7944.1632
$ws.Range("M131").Value = -331324.665
$ws.Range("N131").Value = -18024.1632
$ws.Range("H137").Value = 3017.7778
$ws.Range("I137").Value = 2252
$ws.Range("K137").Value = 6756
$ws.Range("M137").Value = -1656

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 14330.6
$ws.Range("I43").Value = 10408.5
$ws.Range("K43").Value = 10408.5
$ws.Range("M43").Value = -10257.5
$ws.Range("H46").Value = 36999.5
$ws.Range("I46").Value = 36999.5
$ws.Range("K46").Value = 36999.5
$ws.Range("M46").Value = -36843.5
$ws.Range("H57").Value = 13149.25
$ws.Range("H97").Value = 2230.0386
$ws.Range("I97").Value = 2142.0476
$ws.Range("J97").Value = 2599.6
$ws.Range("K97").Value = 2142.0476
$ws.Range("L97").Value = 2599.6
$ws.Range("M97").Value = -1646.0476
$ws.Range("N97").Value = -3591.6
$ws.Range("H102").Value = 2258.4119
$ws.Range("I102").Value = 2055.8125
$ws.Range("K102").Value = 2055.8125
$ws.Range("M102").Value = -433.8125
$ws.Range("H122").Value = 14030.04
$ws.Range("I122").Value = 16355.211
$ws.Range("J122").Value = 6667
$ws.Range("K122").Value = 49065.633
$ws.Range("L122").Value = 20001
$ws.Range("M122").Value = -46615.633
$ws.Range("N122").Value = -24901
$ws.Range("H123").Value = 56569.08
$ws.Range("J123").Value = 56569.08
$ws.Range("L123").Value = 56569.08
$ws.Range("N123").Value = -61469.08
$ws.Range("H132").Value = 14358.568
$ws.Range("I132").Value = 15301.269
$ws.Range("K132").Value = 45903.807
$ws.Range("M132").Value = -43373.807

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3633333
$ws.Range("J2").Value = 4074999.5
$ws.Range("L2").Value = 4074999.5
$ws.Range("N2").Value = -4075223.5
$ws.Range("H16").Value = 2680.276
$ws.Range("I16").Value = 2509.0386
$ws.Range("J16").Value = 4164.3335
$ws.Range("K16").Value = 2509.0386
$ws.Range("L16").Value = 4164.3335
$ws.Range("M16").Value = -2339.0386
$ws.Range("N16").Value = -4504.3335
$ws.Range("H22").Value = 2858.1667
$ws.Range("I22").Value = 2716.3333
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 2716.3333
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -2421.3333
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 2858.1667
$ws.Range("I27").Value = 2716.3333
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 2716.3333
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2609.3333
$ws.Range("N27").Value = -3214
$ws.Range("H61").Value = 3225.5
$ws.Range("I61").Value = 2470.6667
$ws.Range("K61").Value = 2470.6667
$ws.Range("M61").Value = -2268.6667
$ws.Range("H93").Value = 696696.3
$ws.Range("I93").Value = 928222
$ws.Range("J93").Value = 2119.25
$ws.Range("K93").Value = 928222
$ws.Range("L93").Value = 2119.25
$ws.Range("M93").Value = -926974
$ws.Range("N93").Value = -4615.25
$ws.Range("H113").Value = 3225.5
$ws.Range("I113").Value = 2470.6667
$ws.Range("K113").Value = 2470.6667
$ws.Range("M113").Value = -300.6667000000002
$ws.Range("H122").Value = 6923.125
$ws.Range("I122").Value = 6077
$ws.Range("K122").Value = 18231
$ws.Range("M122").Value = -15781
$ws.Range("H136").Value = 5800
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1397.2632
$ws.Range("I107").Value = 931.0833
$ws.Range("K107").Value = 2793.2499
$ws.Range("M107").Value = -873.2498999999998
$ws.Range("H122").Value = 5068.5757
$ws.Range("J122").Value = 27724.75
$ws.Range("L122").Value = 83174.25
$ws.Range("N122").Value = -88074.25
$ws.Range("H126").Value = 3382
$ws.Range("I126").Value = 2917.889
$ws.Range("J126").Value = 6166.6665
$ws.Range("K126").Value = 8753.667000000001
$ws.Range("L126").Value = 18499.9995
$ws.Range("M126").Value = -6283.667000000001
$ws.Range("N126").Value = -23439.9995
$ws.Range("H132").Value = 4088.138
$ws.Range("I132").Value = 4917.4736
$ws.Range("J132").Value = 2512.4
$ws.Range("K132").Value = 14752.4208
$ws.Range("L132").Value = 7537.200000000001
$ws.Range("M132").Value = -12222.4208
$ws.Range("N132").Value = -12597.2
$ws.Range("H133").Value = 64616.668
$ws.Range("J133").Value = 64616.668
$ws.Range("L133").Value = 64616.668
$ws.Range("N133").Value = -74736.66800000001
$ws.Range("H138").Value = 89796.336
$ws.Range("J138").Value = 89695
$ws.Range("L138").Value = 89695
$ws.Range("N138").Value = -99975
